$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 20, pushing existing rows 20.. down to 22..
$ws.Range("A20:R21").Insert()

# Populate the two newly-inserted rows (20 and 21) with their new data
# Row 20
$ws.Cells.Item(20, "A").Value2 = 7
$ws.Cells.Item(20, "B").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(20, "C").Value2 = "Ñuble"
$ws.Cells.Item(20, "D").Value2 = 45222
$ws.Cells.Item(20, "E").Value2 = 16
$ws.Cells.Item(20, "F").Value2 = 100112022
$ws.Cells.Item(20, "G").Value2 = "Arveja Verde"
$ws.Cells.Item(20, "H").Value2 = "Sin especificar"
$ws.Cells.Item(20, "I").Value2 = "Primera"
$ws.Cells.Item(20, "J").Value2 = 20
$ws.Cells.Item(20, "K").Value2 = 25000
$ws.Cells.Item(20, "L").Value2 = 25000
$ws.Cells.Item(20, "M").Value2 = 25000
$ws.Cells.Item(20, "N").Value2 = "`$/saco 25 kilos"
$ws.Cells.Item(20, "O").Value2 = "Región de Ñuble"
$ws.Cells.Item(20, "P").Value2 = 1000
$ws.Cells.Item(20, "Q").Value2 = 25
$ws.Cells.Item(20, "R").Value2 = "Hortaliza"

# Row 21
$ws.Cells.Item(21, "A").Value2 = 7
$ws.Cells.Item(21, "B").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(21, "C").Value2 = "Ñuble"
$ws.Cells.Item(21, "D").Value2 = 45222
$ws.Cells.Item(21, "E").Value2 = 16
$ws.Cells.Item(21, "F").Value2 = 100112022
$ws.Cells.Item(21, "G").Value2 = "Arveja Verde"
$ws.Cells.Item(21, "H").Value2 = "Sin especificar"
$ws.Cells.Item(21, "I").Value2 = "Primera"
$ws.Cells.Item(21, "J").Value2 = 20
$ws.Cells.Item(21, "K").Value2 = 26000
$ws.Cells.Item(21, "L").Value2 = 26000
$ws.Cells.Item(21, "M").Value2 = 26000
$ws.Cells.Item(21, "N").Value2 = "`$/saco 25 kilos"
$ws.Cells.Item(21, "O").Value2 = "Región del Maule"
$ws.Cells.Item(21, "P").Value2 = 1040
$ws.Cells.Item(21, "Q").Value2 = 25
$ws.Cells.Item(21, "R").Value2 = "Hortaliza"

